$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.790.31"
$ws.Range("E2").Value = "  +7.24%  "
$ws.Range("D3").Value = "1.763.23"
$ws.Range("E3").Value = "  +5.69%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'316.68"
$ws.Range("E5").Value = "  +3.09%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "'0.3824"
$ws.Range("E7").Value = "  +3.16%  "
$ws.Range("D8").Value = "'0.3606"
$ws.Range("E8").Value = "  +5.03%  "
$ws.Range("D9").Value = "'50.29"
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "'1.226"
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("E11").Value = "  +6.18%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'21.67"
$ws.Range("E13").Value = "  +6.15%  "
$ws.Range("D14").Value = "'6.468"
$ws.Range("E14").Value = "  +7.45%  "
$ws.Range("D15").Value = "'7.095"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "1.764.03"
$ws.Range("E16").Value = "  +5.82%  "
$ws.Range("D17").Value = "'0.00001154"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").Value = "'0.06792"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("E20").Value = "  +6.33%  "
$ws.Range("D21").Value = "'17.68"
$ws.Range("E21").Value = "  +7.87%  "
$ws.Range("D22").Value = "'6.518"
$ws.Range("E22").Value = "  +6.46%  "
$ws.Range("D23").Value = "'13.04"
$ws.Range("E23").Value = "  +9.07%  "
$ws.Range("D24").Value = "25.755.05"
$ws.Range("E24").Value = "  +7.14%  "
$ws.Range("D25").Value = "'2.431"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("D26").Value = "'2.905"
$ws.Range("E26").Value = "  +9.19%  "
$ws.Range("D27").Value = "'20.77"
$ws.Range("D28").Value = "'155.97"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").Value = "1.961.08"
$ws.Range("E29").Value = "  +6.03%  "
$ws.Range("D30").Value = "'133.80"
$ws.Range("E30").Value = "  +5.33%  "
$ws.Range("D31").Value = "'1.206"
$ws.Range("E31").Value = "  +23.27%  "
$ws.Range("D32").Value = "'7.189"
$ws.Range("E32").Value = "  +13.19%  "
$ws.Range("D33").Value = "'4.219"
$ws.Range("E33").Value = "  +3.84%  "
$ws.Range("D34").Value = "'14.36"
$ws.Range("E34").Value = "  +16.99%  "
$ws.Range("D35").Value = "'1.808"
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("D36").Value = "'0.08768"
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("D37").Value = "'5.721"
$ws.Range("E37").Value = "  +7.76%  "

# Row 38/39: Hedera and VeChain swap places
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02503"
$ws.Range("E38").Value = "  +8.14%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06751"
$ws.Range("E39").Value = "  +6.42%  "
$ws.Range("D40").Value = "'9.386"
$ws.Range("E40").Value = "  +5.26%  "
$ws.Range("D41").Value = "'0.2261"
$ws.Range("E41").Value = "  +9.03%  "
$ws.Range("D42").Value = "'1.297"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'0.6589"
$ws.Range("E43").Value = "  +8.34%  "
$ws.Range("D44").Value = "'14.36"
$ws.Range("E44").Value = "  +8.86%  "
$ws.Range("D45").Value = "'0.9995"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "'0.6350"
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("D47").Value = "'3.907"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("E48").Value = "  +8.66%  "
$ws.Range("D49").Value = "'132.06"
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("D50").Value = "'0.07509"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("D51").Value = "'81.14"
$ws.Range("E51").Value = "  +6.98%  "
